$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add I0 and IF headers, matching the style of the existing
#     header cells (H1 etc., style index 1: bold font + border + centered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-34: I column is always 1 ("I0"), J column mirrors H
#     column's existing "IP" value (effectively "IF" = IP here).
for ($row = 2; $row -le 34; $row++) {
    $hValue = $ws.Cells.Item($row, 8).Value2
    $ws.Cells.Item($row, 9).Value = 1
    $ws.Cells.Item($row, 10).Value = $hValue
}
